# Generate Report for Handoff
# Updates the localization-status workbook to reflect a newly generated
# handoff report: refreshes the "Latest Handoff Datetime" timestamps for
# the zh-cn and de-de sheets (and the corresponding "Latest HO Xliff
# Generate Date" roll-up on the Overview sheet), and marks the handoff
# Priority ("ht") for the rows that were just handed off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 11, 13, 14)

# zh-cn: "Latest Handoff Datetime" (column H) moves forward to 16:23:06
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-23 16:23:06"
}

# de-de: "Latest Handoff Datetime" (column H) moves forward to 16:23:18
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-23 16:23:18"
}

# Overview: "Latest HO Xliff Generate Date" (column G) reflects the newest
# handoff timestamp across locales, which is now 16:23:18
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-23 16:23:18"
}

# Mark Priority ("ht" = hot/handoff) for the files that were just handed off
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
